$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Re-enter the radius.food (column E) formulas as fill-down ranges so Excel
# stores them as shared formulas (t="shared"), matching the re-saved file ---
$ws.Range("E2:E33").Formula = "=D2/2"
$ws.Range("E34:E65").Formula = "=D34/2"
$ws.Range("E66:E97").Formula = "=D66/2"
# The third fill-down extended one group past the last data row (92); drop the
# spill-over cells again so no phantom rows 93:97 are introduced (the shared
# formula's "ref" span on E66 still records the original E66:E97 fill extent).
$ws.Range("E93:E97").ClearContents()

# --- Row 73: corrected radius.galea reading (C73) ---
$ws.Range("C73").Value = 0.15402083333333333

# --- Rows 75/76: the two diam.food (D) readings were swapped ---
$ws.Range("D75").Value = 0.21871200000000002
$ws.Range("D76").Value = 0.23214533333333334

# --- Rows 77/78: the two radius.galea (C) readings were swapped ---
$ws.Range("C77").Value = 0.08806249999999999
$ws.Range("C78").Value = 0.10047222222222223

# --- Rows 86-92 (Manduca_quinquemaculata block): corrected readings ---
$ws.Range("C86").Value = 0.24597916666666669
$ws.Range("D86").Value = 0.22288550000000001

$ws.Range("C87").Value = 0.18454166666666666
$ws.Range("D87").Value = 0.21437299999999998

$ws.Range("C88").Value = 0.14569444444444438
$ws.Range("D88").Value = 0.20825199999999999

$ws.Range("C89").Value = 0.13672916666666668
$ws.Range("D89").Value = 0.19321099999999999

$ws.Range("C90").Value = 0.13966666666666666
$ws.Range("D90").Value = 0.18458350000000001

$ws.Range("C91").Value = 0.12945833333333334
$ws.Range("D91").Value = 0.17958350000000001

$ws.Range("C92").Value = 0.10677777777777779
$ws.Range("D92").Value = 0.16642266666666666

# --- Update the active selection to match the final saved view ---
[void]$ws.Range("H86:I86").Select()
